$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "50.892.50"
$ws.Range("E2").Value = "  -1.77%  "
$ws.Range("D3").Value = "2.741.07"
$ws.Range("E3").Value = "  -1.28%  "
$ws.Range("D4").Value = "'0.999"
$ws.Range("E4").Value = "  -0.06%  "
$ws.Range("D5").Value = "'350.39"
$ws.Range("E5").Value = "  -2.23%  "
$ws.Range("D6").Value = "'106.79"
$ws.Range("E6").Value = "  -1.94%  "
$ws.Range("D7").Value = "'0.544"
$ws.Range("E7").Value = "  -2.39%  "
$ws.Range("E8").Value = "  +0.02%  "
$ws.Range("D9").Value = "'0.575"
$ws.Range("E9").Value = "  -2.50%  "
$ws.Range("D10").Value = "'38.71"
$ws.Range("E10").Value = "  -2.98%  "
$ws.Range("E11").Value = "  +3.58%  "
$ws.Range("D12").Value = "'0.0827"
$ws.Range("E12").Value = "  -2.29%  "
$ws.Range("D13").Value = "'19.43"
$ws.Range("E13").Value = "  +0.19%  "
$ws.Range("D14").Value = "'7.40"
$ws.Range("E14").Value = "  -2.80%  "
$ws.Range("D15").Value = "3.164.95"
$ws.Range("E15").Value = "  -1.52%  "
$ws.Range("D16").Value = "2.728.58"
$ws.Range("E16").Value = "  -2.17%  "
$ws.Range("D17").Value = "'0.912"
$ws.Range("E17").Value = "  -0.06%  "
$ws.Range("D18").Value = "50.810.00"
$ws.Range("E18").Value = "  -1.59%  "
$ws.Range("E19").Value = "  +2.49%  "
$ws.Range("D20").Value = "'3.01"
$ws.Range("E20").Value = "  -2.90%  "
$ws.Range("D21").Value = "'12.87"
$ws.Range("E21").Value = "  -1.45%  "
$ws.Range("E22").Value = "  -2.91%  "
$ws.Range("D23").Value = "'68.79"
$ws.Range("E23").Value = "  -0.86%  "
$ws.Range("D24").Value = "'261.87"
$ws.Range("E24").Value = "  -4.28%  "
$ws.Range("E25").Value = "  -2.19%  "
$ws.Range("E26").Value = "  +0.10%  "
$ws.Range("D27").Value = "'25.70"
$ws.Range("E27").Value = "  -2.57%  "
$ws.Range("D28").Value = "'0.159"
$ws.Range("E28").Value = "  +12.34%  "
$ws.Range("E29").Value = "  +0.12%  "
$ws.Range("E30").Value = "  -1.57%  "
$ws.Range("D31").Value = "'51.64"
$ws.Range("E31").Value = "  +0.16%  "
$ws.Range("D32").Value = "'34.18"
$ws.Range("E32").Value = "  +0.12%  "
$ws.Range("D33").Value = "'5.93"
$ws.Range("E33").Value = "  +3.82%  "
$ws.Range("D34").Value = "'0.0436"
$ws.Range("E34").Value = "  -7.12%  "
$ws.Range("E35").Value = "  -1.91%  "
$ws.Range("D36").Value = "'5.14"
$ws.Range("E36").Value = "  -4.42%  "
$ws.Range("D37").Value = "'0.999"
$ws.Range("E37").Value = "  -0.05%  "
$ws.Range("D38").Value = "'18.36"
$ws.Range("E38").Value = "  +2.26%  "
$ws.Range("E39").Value = "  -1.59%  "
$ws.Range("E40").Value = "  -3.13%  "
$ws.Range("E41").Value = "  -1.71%  "
$ws.Range("D42").Value = "'2.45"
$ws.Range("E42").Value = "  -2.31%  "
$ws.Range("D43").Value = "'120.54"
$ws.Range("E43").Value = "  -3.90%  "
$ws.Range("E44").Value = "  -2.44%  "
$ws.Range("D45").Value = "'21.66"
$ws.Range("E45").Value = "  -0.62%  "
$ws.Range("D46").Value = "2.071.54"
$ws.Range("E46").Value = "  +0.97%  "
$ws.Range("D47").Value = "'2.28"
$ws.Range("E47").Value = "  -1.90%  "
$ws.Range("E48").Value = "  -1.12%  "
$ws.Range("E49").Value = "  -2.16%  "
$ws.Range("D50").Value = "'5.40"
$ws.Range("E50").Value = "  -5.58%  "
$ws.Range("E51").Value = "  +5.25%  "
